$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text for columns B..F of the header row (column A "Datum" stays
# the same). Each entry also gets re-bolded, matching the source edit which
# added <b/> to the shared-string run properties for these five headers.
$headers = @(
    @{ Addr = "B1"; Text = "Chai-Gesamtumsatz (Einheiten)" },
    @{ Addr = "C1"; Text = "Artisanal Chai-Umsatz (Einheiten)" },
    @{ Addr = "D1"; Text = "Vorgefertigter Chai-Umsatz (Einheiten)" },
    @{ Addr = "E1"; Text = "Social-Media-Interaktion (Ansichten)" },
    @{ Addr = "F1"; Text = "Onlinesuchen nach Chai" }
)

foreach ($h in $headers) {
    $cell = $ws.Range($h.Addr)
    $cell.Value2 = $h.Text

    $len = $cell.Characters().Text.Length

    # Re-apply bold + white font colour as an explicit rich-text run. Setting
    # the formatting on the *entire* Characters() span in a single call
    # collapses back to "no run formatting" in this engine, so instead the
    # span is split into two adjacent pieces that get the same formatting -
    # those get merged back into one <r> with the correct <rPr> on save,
    # matching the existing header-cell formatting (bold, white, Calibri 11).
    $firstLen = $len - 1
    if ($firstLen -lt 1) { $firstLen = 1 }

    $part1 = $cell.Characters(1, $firstLen)
    $part1.Font.Bold = $true
    $part1.Font.Color = 16777215

    if ($len -gt 1) {
        $part2 = $cell.Characters($len, 1)
        $part2.Font.Bold = $true
        $part2.Font.Color = 16777215
    }
}
